# Weekly update: insert a new daily record as row 14 (shifting the
# existing rows 14-82 down to 15-83) and fill in its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; Excel shifts rows 14:82 down to 15:83
# and the sheet dimension grows from R82 to R83 automatically.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(14, 3).Value = "Maule"
$ws.Cells.Item(14, 4).Value = 44881
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = 100112040
$ws.Cells.Item(14, 7).Value = "Cilantro"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 150
$ws.Cells.Item(14, 11).Value = 7000
$ws.Cells.Item(14, 12).Value = 7000
$ws.Cells.Item(14, 13).Value = 7000
$ws.Cells.Item(14, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 194
$ws.Cells.Item(14, 17).Value = 36
$ws.Cells.Item(14, 18).Value = "Hortaliza"
